$d = $word.ActiveDocument

# 1) Remove the trailing space at the end of the "Com esse gráfico..." paragraph.
$d.Content.Find.Execute(
    "entre 3 e 4. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "entre 3 e 4.", 2
)

# 2) Append the new list paragraphs after the last paragraph in the document.
$lastP = $d.Paragraphs.Last
$lastR = $lastP.Range
$lastR.InsertParagraphAfter()

# --- Paragraph: " O que aprendemos:" (ilvl = 1, numId = 2) ---
$p1 = $d.Paragraphs.Last
$p1.Range.Text = " O que aprendemos:"
$p1.Range.ListFormat.ListLevelNumber = 2
$p1.Range.InsertParagraphAfter()

# --- Paragraph: "Importar pandas que é uma biblioteca de análise de dados;" (ilvl = 2) ---
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Importar pandas que é uma biblioteca de análise de dados;"
$p2.Range.ListFormat.ListLevelNumber = 3
$p2.Range.InsertParagraphAfter()

# --- Paragraph: "Ler dados CSV;" (ilvl = 2) ---
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "Ler dados CSV;"
$p3.Range.ListFormat.ListLevelNumber = 3
$p3.Range.InsertParagraphAfter()

# --- Paragraph: "Fazer upload de arquivo para o Google Colab;" (ilvl = 2) ---
$p4 = $d.Paragraphs.Last
$p4.Range.Text = "Fazer upload de arquivo para o Google Colab;"
$p4.Range.ListFormat.ListLevelNumber = 3
$p4.Range.InsertParagraphAfter()

# --- Paragraph: "Renomear colunas;" (ilvl = 2) ---
$p5 = $d.Paragraphs.Last
$p5.Range.Text = "Renomear colunas;"
$p5.Range.ListFormat.ListLevelNumber = 3
$p5.Range.InsertParagraphAfter()

# --- Paragraph: "Contar dados;" (ilvl = 2) ---
$p6 = $d.Paragraphs.Last
$p6.Range.Text = "Contar dados;"
$p6.Range.ListFormat.ListLevelNumber = 3
$p6.Range.InsertParagraphAfter()

# --- Paragraph: "Melhorar a visualização das informações. " (ilvl = 2), two runs ---
$p7 = $d.Paragraphs.Last
$p7.Range.Text = "Melhorar a visualização das informações."
$p7.Range.ListFormat.ListLevelNumber = 3
$p7.Range.InsertAfter(" ")

Write-Output "done"
